# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
#
# The Price column stores numeric-looking values (e.g. "30.285.91",
# "0.9998") as plain TEXT in the workbook (t="inlineStr"), not as real
# numbers. Excel's Range.Value setter auto-detects numeric-looking
# strings and coerces them to the Number type (which would also silently
# drop meaningful trailing zeros, e.g. "255.80" -> 255.8). To faithfully
# reproduce the original text values we temporarily force the cell's
# NumberFormat to Text ("@") before assigning the value, then clear the
# format again afterwards so the cell's style index is left exactly as
# it was (General/default) — only the stored value changes.
#
# The Volume(1h) column values (e.g. "  -0.21%  ") are never numeric
# (percent sign + padding spaces), so they can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.285.91"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.21%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.930.83"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -0.32%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "248.95"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.33%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.7162"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -1.11%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.02%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3195"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -4.31%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "27.71"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -2.97%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07109"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -4.70%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.7917"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -2.89%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07984"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.79%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.924.13"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -2.23%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "94.86"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  -1.11%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "30.277.42"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -0.19%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "255.80"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.39%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008040"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -4.04%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "5.779"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -1.58%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.179.97"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -0.35%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9994"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -0.01%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.04%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.822"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  -2.74%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "165.12"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +2.90%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "19.11"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -1.78%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.267"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -6.30%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.1270"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -4.83%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.354"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +0.84%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.527"
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.393"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -1.15%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.133"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  -1.13%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.271"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +1.93%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7454"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -0.71%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.774"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +1.49%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01960"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -1.82%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.799"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "78.61"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.02%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "6.356"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -4.71%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.4509"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.59%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.987"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.89%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8471"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +0.85%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.9995"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.10%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "100.49"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("E47").Value = "  -0.33%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.417"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.41%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "36.67"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -0.43%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "944.97"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +9.75%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.06105"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +2.07%  "
